$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 1000000000
$ws.Range("I5").Value = 1000000000
$ws.Range("K5").Value = 1000000000
$ws.Range("M5").Value = -999999885
$ws.Range("H41").Value = 199.25
$ws.Range("I41").Value = 199.25
$ws.Range("K41").Value = 199.25
$ws.Range("M41").Value = 240.75
$ws.Range("H43").Value = 11498.5
$ws.Range("I43").Value = 11498.5
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 11498.5
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -11429.5
$ws.Range("H105").Value = 46937.11
$ws.Range("I105").Value = 60000
$ws.Range("J105").Value = 43204.855
$ws.Range("K105").Value = 60000
$ws.Range("L105").Value = 43204.855
$ws.Range("M105").Value = -56506
$ws.Range("N105").Value = -50192.855
$ws.Range("H125").Value = 993.6667
$ws.Range("I125").Value = 997.5
$ws.Range("K125").Value = 8977.5
$ws.Range("M125").Value = -6517.5
$ws.Range("H138").Value = 3646.5278
$ws.Range("J138").Value = 5894.7144
$ws.Range("L138").Value = 17684.1432
$ws.Range("N138").Value = -27964.1432

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 860
$ws.Range("I22").Value = 790
$ws.Range("K22").Value = 790
$ws.Range("M22").Value = -491
$ws.Range("H32").Value = 8910.817999999999
$ws.Range("I32").Value = 6454.185
$ws.Range("K32").Value = 6454.185
$ws.Range("M32").Value = -6167.185
$ws.Range("H41").Value = 4230.769
$ws.Range("I41").Value = 4230.769
$ws.Range("K41").Value = 4230.769
$ws.Range("M41").Value = -3816.769
$ws.Range("H94").Value = 37665
$ws.Range("J94").Value = 37665
$ws.Range("L94").Value = 37665
$ws.Range("N94").Value = -39467

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 2736.25
$ws.Range("I36").Value = 2315
$ws.Range("J36").Value = 4000
$ws.Range("K36").Value = 2315
$ws.Range("L36").Value = 4000
$ws.Range("M36").Value = -1781
$ws.Range("N36").Value = -5068

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3629.818
$ws.Range("I31").Value = 2525
$ws.Range("K31").Value = 2525
$ws.Range("M31").Value = -2230
$ws.Range("H34").Value = 3629.818
$ws.Range("I34").Value = 2525
$ws.Range("K34").Value = 2525
$ws.Range("M34").Value = -2323
$ws.Range("H88").Value = 35239.547
$ws.Range("J88").Value = 35239.547
$ws.Range("L88").Value = 35239.547
$ws.Range("N88").Value = -36051.547
$ws.Range("H91").Value = 35239.547
$ws.Range("J91").Value = 35239.547
$ws.Range("L91").Value = 35239.547
$ws.Range("N91").Value = -38047.547

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 59293.332
$ws.Range("J93").Value = 59293.332
$ws.Range("L93").Value = 59293.332
$ws.Range("N93").Value = -63037.332
$ws.Range("H98").Value = 18237
$ws.Range("J98").Value = 18237
$ws.Range("L98").Value = 18237
$ws.Range("N98").Value = -24227
$ws.Range("H113").Value = 3812.25
$ws.Range("I113").Value = 3666.3333
$ws.Range("K113").Value = 3666.3333
$ws.Range("M113").Value = -1496.3333
$ws.Range("H126").Value = 5763
$ws.Range("I126").Value = 5512
$ws.Range("K126").Value = 16536
$ws.Range("M126").Value = -14066

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2250
$ws.Range("I7").Value = 2250
$ws.Range("K7").Value = 2250
$ws.Range("M7").Value = -2138
$ws.Range("H16").Value = 1411.6111
$ws.Range("I16").Value = 1254.8334
$ws.Range("J16").Value = 1725.1666
$ws.Range("K16").Value = 1254.8334
$ws.Range("L16").Value = 1725.1666
$ws.Range("M16").Value = -1084.8334
$ws.Range("N16").Value = -2065.1666
$ws.Range("H22").Value = 2987.5334
$ws.Range("I22").Value = 1210.909
$ws.Range("K22").Value = 1210.909
$ws.Range("M22").Value = -915.9090000000001
$ws.Range("H27").Value = 2987.5334
$ws.Range("I27").Value = 1210.909
$ws.Range("K27").Value = 1210.909
$ws.Range("M27").Value = -1103.909
$ws.Range("H33").Value = 15000000
$ws.Range("I33").Value = 15000000
$ws.Range("K33").Value = 15000000
$ws.Range("M33").Value = -14999710
$ws.Range("H40").Value = 3099.4
$ws.Range("I40").Value = 3099.4
$ws.Range("K40").Value = 3099.4
$ws.Range("M40").Value = -2963.4
$ws.Range("H55").Value = 807.6667
$ws.Range("I55").Value = 500
$ws.Range("J55").Value = 1192.25
$ws.Range("K55").Value = 500
$ws.Range("L55").Value = 1192.25
$ws.Range("M55").Value = -327
$ws.Range("N55").Value = -1538.25
$ws.Range("H68").Value = 3050.75
$ws.Range("I68").Value = 2400
$ws.Range("K68").Value = 2400
$ws.Range("M68").Value = -1651
$ws.Range("H71").Value = 3050.75
$ws.Range("I71").Value = 2400
$ws.Range("K71").Value = 12000
$ws.Range("M71").Value = -8256
$ws.Range("H82").Value = 3607.4
$ws.Range("I82").Value = 4012.3333
$ws.Range("J82").Value = 3000
$ws.Range("K82").Value = 4012.3333
$ws.Range("L82").Value = 3000
$ws.Range("M82").Value = -3651.3333
$ws.Range("N82").Value = -3722
$ws.Range("H85").Value = 3607.4
$ws.Range("I85").Value = 4012.3333
$ws.Range("J85").Value = 3000
$ws.Range("K85").Value = 4012.3333
$ws.Range("L85").Value = 3000
$ws.Range("M85").Value = -2764.3333
$ws.Range("N85").Value = -5496
$ws.Range("H126").Value = 2250
$ws.Range("I126").Value = 2250
$ws.Range("K126").Value = 6750
$ws.Range("M126").Value = -4280
$ws.Range("H132").Value = 2500
$ws.Range("I132").Value = 2050.5334
$ws.Range("J132").Value = 5196.8
$ws.Range("K132").Value = 6151.600199999999
$ws.Range("L132").Value = 15590.4
$ws.Range("M132").Value = -3621.600199999999
$ws.Range("N132").Value = -20650.4

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H122").Value = 2200.3845
$ws.Range("J122").Value = 1960.5
$ws.Range("L122").Value = 5881.5
$ws.Range("N122").Value = -10781.5
$ws.Range("H126").Value = 3077.3333
$ws.Range("J126").Value = 3399.6667
$ws.Range("L126").Value = 10199.0001
$ws.Range("N126").Value = -15139.0001
